$wb = $excel.ActiveWorkbook

# Duplicate the active "Sedan_HambaLG_f" sheet (the template for a single
# Droplink hardpoint set) to create a new sheet for the Trailer1Axle_f
# front droplink, placed after the last existing sheet.
$template = $wb.Worksheets.Item("Sedan_HambaLG_f")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Trailer1Axle_f"

# Update the Instance label for the new droplink.
$newSheet.Range("H3").Value = "Droplink_Trailer1Axle_f"

# Updated hardpoints / values for the Trailer1Axle_f droplink (v2p14,
# compatible with MF-Swift v2212).
$newSheet.Range("F5").Value = 0.05
$newSheet.Range("G5").Value = 0.6
$newSheet.Range("H5").Value = 0.19

$newSheet.Range("F6").Formula = "=0.3-0.15"
$newSheet.Range("G6").Value = 0.57999999999999996
$newSheet.Range("H6").Value = 0.2

$newSheet.Range("H7").Value = 50

$newSheet.Range("H7").Select()
